# Updated symbol list on Mon Feb  6 14:25:50 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) text cells for
# the coin rows whose market data changed since the previous snapshot.
# The sheet stores these figures as plain text (e.g. "327.48", "-0.58%"),
# so each write is forced to text with a leading apostrophe; the style is
# then reset to "Normal" because Excel stamps a quote-prefix number
# format on text that merely looks numeric, which the source file never had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.58%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'43.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.56%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.556"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.87%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'-1.91%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.900"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.09%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.275"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.05%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E9").Value = "'0.08%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-3.02%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1840"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-3.99%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09652"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.18%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.04391"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.79%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1067"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.19%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001287"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.40%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005984"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.29%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.403"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.99%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.3491"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.27%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'9.996"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'14.43%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.78%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2507"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.59%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04202"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.59%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001248"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.57%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004286"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.36%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001262"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.12%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003994"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.32%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02633"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-7.04%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05471"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-4.62%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007568"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-4.41%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1392"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.64%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.008009"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-18.24%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-4.28%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008834"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-12.41%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006932"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-5.28%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.33%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.002273"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.32%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003561"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1.76%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.33%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.33%"
$ws.Range("E50").Style = "Normal"
